$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8   = -12.10199999999999
    10  = -13.6011
    12  = -10.3856
    18  = -12.0197
    25  = -11.6856
    37  = -13.1118
    55  = -13.52189999999999
    68  = -11.5615
    77  = -13.92010000000001
    78  = -13.8653
    79  = -13.7024
    80  = -13.76000000000001
    81  = -14.3435
    82  = -12.2427
    84  = -13.77829999999999
    101 = -12.5759
    102 = -12.098
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
